# Apply roster update: remove "Nicolas Claxton" row, re-order remaining
# player rows so "Anthony Edwards" sits right after "Jose Alvarado", and
# shrink the sheet dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired table (header stays in row 1; data rows 2-17)
$data = @(
    @("Jose Alvarado", "PG", "New Orleans Pelicans"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Amen Thompson", "SG,SF,PF", "Houston Rockets"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers"),
    @("RJ Barrett", "SG,SF,PF", "Toronto Raptors"),
    @("Draymond Green", "PF,C", "Golden State Warriors"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Fred VanVleet", "PG", "Houston Rockets")
)

# Clear out the old data rows (2 through 18) before writing the new,
# shorter table back so no stale cells remain below row 17.
$oldLastRow = 18
$clearRange = $ws.Range("A2:C$oldLastRow")
$clearRange.Clear() | Out-Null

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

$wb.Save()
